# Applies the cryptos price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.594.91"
$ws.Range("E2").Value = "  -2.41%  "
$ws.Range("D3").Value = "2.895.01"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.04%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "2.892.73"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D16").Value = "3.374.89"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").Value = "61.610.67"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "2.894.07"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.657"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.51%  "
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("E32").Value = "  -7.24%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.61%  "
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.06%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("E44").Value = "  -5.24%  "
$ws.Range("D45").Value = "2.691.94"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0336"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "348.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.16%  "
